$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the "Test case 1" block (rows 1,3-9) below as "Test case 2" (rows 11,13-19),
# copying formatting along with values/merges.
$ws.Range("A1:F1").Copy($ws.Range("A11:F11"))
$ws.Range("A3:D4").Copy($ws.Range("A13:D14"))
$ws.Range("A5:D5").Copy($ws.Range("A15:D15"))
$ws.Range("A7:F7").Copy($ws.Range("A17:F17"))
$ws.Range("A8:F9").Copy($ws.Range("A18:F19"))

# Row 1's bigger title font drives a taller row; match that on the new title row too.
$ws.Rows(11).RowHeight = $ws.Rows(1).RowHeight

# Update the text that differs between the two test cases.
$ws.Range("A11").Value2 = "Test case 2"
$ws.Range("B13").Value2 = "UT_002"

# The second test case exercises different Test Data / Expectations values.
$ws.Range("D18").Value2 = 2
$ws.Range("E18").Value2 = 2
$ws.Range("C19").Value2 = 3

# Reflect the view position/selection saved with the workbook after the edit.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("B22").Select() | Out-Null

Write-Output "done"
